# Delete p3 (row 4) because IMU recording was only done for the second ascent.
#
# The physiotope labels in column A (p1, p2, p3, ...) are a fixed sequential
# listing that is left untouched for the rows that remain; only the
# measurement columns (B:I) for every physiotope after p3 move up one row,
# and the final row (which held p17's measurements) becomes blank and drops
# out of the used range.
#
# This engine's Range.Delete(xlShiftUp) shifts every column of the row
# (including column A), so we delete the whole row first and then restore
# the column-A labels for the rows that are left (p3..p16), row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("4:4").Delete()

for ($r = 4; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = "p" + ($r - 1)
}

$ws.Range("A19").Select()
